# The post at row 667 ("「惑星地球はあなた方を歓迎します」...") was removed.
# Deleting the entire row shifts all subsequent rows (668-830) up by one,
# which matches the target diff (dimension shrinks from A1:C830 to A1:C829).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(667).Delete()
